$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 87
$ws.Range("A87").Value = 229
$ws.Range("C87").Value = 84
$ws.Range("F87").Value = 87
$ws.Range("G87").Value = 88
$ws.Range("H87").Value = 8

# Row 88
$ws.Range("F88").Value = 88
$ws.Range("H88").Value = 0

# Row 89
$ws.Range("A89").Value = 237
$ws.Range("B89").Value = 156
$ws.Range("C89").Value = 74
$ws.Range("F89").Value = 89
$ws.Range("G89").Value = 78
$ws.Range("H89").Value = 8

# Row 90
$ws.Range("A90").Value = 248
$ws.Range("B90").Value = 163
$ws.Range("C90").Value = 79
$ws.Range("D90").Value = 3
$ws.Range("F90").Value = 90
$ws.Range("G90").Value = 82
$ws.Range("H90").Value = 11

# Row 91
$ws.Range("A91").Value = 249
$ws.Range("B91").Value = 164
$ws.Range("H91").Value = 1
